# ES6 Destructuring doc review edit
# 1. Remove the stray _GoBack bookmark after "... and store in a ".
# 2. Highlight (yellow) several runs across a few paragraphs discussing
#    object-destructured function parameters.
# 3. Split the final "default values if none are provided." sentence and
#    plant a fresh _GoBack bookmark in the middle of it (mirrors where the
#    cursor was left when the author stopped editing).

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark -------------------------------
$d.Bookmarks.Item("_GoBack").Delete()

# Helper: walk forward through a paragraph's range, highlighting each
# text segment (in order) with yellow. Returns nothing; advances through
# the paragraph so repeated/short substrings (e.g. "x") resolve to the
# correct occurrence.
function Highlight-Segments($paragraph, [string[]]$segments) {
    $paraEnd = $paragraph.Range.End
    $cursor = $paragraph.Range.Start
    foreach ($seg in $segments) {
        $r = $d.Range($cursor, $paraEnd)
        $r.Find.ClearFormatting()
        $r.Find.Replacement.ClearFormatting()
        $r.Find.Replacement.Highlight = $true
        [void]$r.Find.Execute($seg, $true, $false, $false, $false, $false, $true, 0, $true, $seg, 1)
        $cursor = $r.End
    }
}

# --- 2a. "Now we can define the function parameter list as an object
#          destructure pattern, like so:" -------------------------------
$p1 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Now we can define the function parameter list*") {
        $p1 = $para
        break
    }
}
Highlight-Segments $p1 @(
    "Now we can define the function parameter list as an ",
    "object ",
    "destructure",
    " pattern"
)

# --- 2b. "Notice that in the function body above we can refer to x
#          directly, ... like options.x." --------------------------------
$p2 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Notice that in the function body above*") {
        $p2 = $para
        break
    }
}
$rightQuote = [char]0x2019
Highlight-Segments $p2 @(
    "Notice that in the function body above we can refer to ",
    "x",
    (" directly, we don" + $rightQuote + "t have to refer to it through an object property like "),
    "options.x"
)

# --- 2c. "In addition to that when using destructured function
#          parameters we can also provide default values, like so:" -----
$p3 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*In addition to that when using*") {
        $p3 = $para
        break
    }
}
Highlight-Segments $p3 @(
    "In addition to that when using ",
    "destructured",
    " function parameters we can also provide default values"
)

# --- 3. Split "... default values if none are provided." and drop a
#        fresh _GoBack bookmark between "value" and "s if none ..." -----
$full = $d.Content.Text
$needle = " syntax for providing optional parameters to functions, including giving them default value"
$idx = $full.IndexOf($needle)
$splitPoint = $idx + $needle.Length
$bmRange = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
